$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
    # MatchCase:=True so only the exact-case typo is touched; Replace:=2 (wdReplaceAll)
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "proudcing" "producing"
Replace-Text "slow down times" "slowdown times"
Replace-Text "Poission distribution" "Poisson distribution"
Replace-Text "analzying parameter" "analyzing parameter"
Replace-Text "quantative finance" "quantitative finance"
Replace-Text "probabalistic distribution" "probabilistic distribution"
Replace-Text "kappa, or the rate parameter for the exponential" "kappa,  or the rate parameter for the exponential"
